$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 362; existing rows 362-376 shift down to 363-377.
$ws.Rows.Item(362).Insert()

# New row 362 (brand-new data point).
$ws.Range("A362").Value = 8
$ws.Range("B362").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C362").Value = 'Coquimbo'
$ws.Range("D362").Value = 45041
$ws.Range("E362").Value = 4
$ws.Range("F362").Value = 100112031
$ws.Range("G362").Value = 'Poroto verde'
$ws.Range("H362").Value = 'Magnum'
$ws.Range("I362").Value = 'Primera'
$ws.Range("J362").Value = 400
$ws.Range("K362").Value = 29000
$ws.Range("L362").Value = 30000
$ws.Range("M362").Value = 29500
$ws.Range("N362").Value = '$/malla 25 kilos'
$ws.Range("O362").Value = 'Provincia de Limarí'
$ws.Range("P362").Value = 1180
$ws.Range("Q362").Value = 25
$ws.Range("R362").Value = 'Hortaliza'

# Update the shifted rows (363-377) where values differ from the source rows they were shifted from.
# Row 363
$ws.Range("D363").Value = 45041
$ws.Range("H363").Value = 'Sin especificar'
$ws.Range("J363").Value = 360
$ws.Range("K363").Value = 30000
$ws.Range("L363").Value = 31000
$ws.Range("M363").Value = 30500
$ws.Range("P363").Value = 1220
# Row 364
$ws.Range("D364").Value = 44705
$ws.Range("K364").Value = 26000
$ws.Range("L364").Value = 27000
$ws.Range("M364").Value = 26500
$ws.Range("P364").Value = 1060
# Row 365
$ws.Range("D365").Value = 44342
$ws.Range("J365").Value = 500
$ws.Range("K365").Value = 30000
$ws.Range("L365").Value = 31000
$ws.Range("M365").Value = 30500
$ws.Range("O365").Value = 'Provincia de Limarí'
$ws.Range("P365").Value = 1220
# Row 366
$ws.Range("D366").Value = 44179
$ws.Range("J366").Value = 400
$ws.Range("K366").Value = 16000
$ws.Range("L366").Value = 17000
$ws.Range("M366").Value = 16500
$ws.Range("O366").Value = 'Provincia del Elquí'
$ws.Range("P366").Value = 660
# Row 367
$ws.Range("D367").Value = 44172
$ws.Range("J367").Value = 600
$ws.Range("K367").Value = 15000
$ws.Range("L367").Value = 16000
$ws.Range("M367").Value = 15500
$ws.Range("O367").Value = 'Provincia de Limarí'
$ws.Range("P367").Value = 620
# Row 368
$ws.Range("D368").Value = 44952
$ws.Range("J368").Value = 400
$ws.Range("K368").Value = 23000
$ws.Range("L368").Value = 24000
$ws.Range("M368").Value = 23500
$ws.Range("O368").Value = 'Provincia del Elquí'
$ws.Range("P368").Value = 940
# Row 369
$ws.Range("D369").Value = 44727
$ws.Range("J369").Value = 480
$ws.Range("K369").Value = 25000
$ws.Range("L369").Value = 26000
$ws.Range("M369").Value = 25500
$ws.Range("P369").Value = 1020
# Row 370
$ws.Range("D370").Value = 44391
$ws.Range("J370").Value = 600
$ws.Range("K370").Value = 22000
$ws.Range("L370").Value = 23000
$ws.Range("M370").Value = 22500
$ws.Range("O370").Value = 'Perú'
$ws.Range("P370").Value = 900
# Row 371
$ws.Range("D371").Value = 44168
$ws.Range("J371").Value = 480
$ws.Range("K371").Value = 19000
$ws.Range("L371").Value = 20000
$ws.Range("M371").Value = 19500
$ws.Range("O371").Value = 'Provincia de Limarí'
$ws.Range("P371").Value = 780
# Row 372
$ws.Range("D372").Value = 44875
$ws.Range("J372").Value = 500
$ws.Range("K372").Value = 41000
$ws.Range("L372").Value = 42000
$ws.Range("M372").Value = 41500
$ws.Range("O372").Value = 'Región de Arica y Parinacota'
$ws.Range("P372").Value = 1660
# Row 373
$ws.Range("D373").Value = 44454
$ws.Range("K373").Value = 33000
$ws.Range("L373").Value = 34000
$ws.Range("M373").Value = 33500
$ws.Range("P373").Value = 1340
# Row 374
$ws.Range("D374").Value = 44426
$ws.Range("J374").Value = 700
$ws.Range("K374").Value = 32500
$ws.Range("L374").Value = 33000
$ws.Range("M374").Value = 32750
$ws.Range("O374").Value = 'Perú'
$ws.Range("P374").Value = 1310
# Row 375
$ws.Range("D375").Value = 44526
$ws.Range("J375").Value = 520
$ws.Range("K375").Value = 22000
$ws.Range("L375").Value = 23000
$ws.Range("M375").Value = 22500
$ws.Range("P375").Value = 900
# Row 376
$ws.Range("D376").Value = 44918
$ws.Range("H376").Value = 'Magnum'
$ws.Range("J376").Value = 400
$ws.Range("O376").Value = 'Provincia de Limarí'
# Row 377
$ws.Range("D377").Value = 44217
$ws.Range("H377").Value = 'Sin especificar'
$ws.Range("J377").Value = 500
$ws.Range("K377").Value = 23000
$ws.Range("M377").Value = 23500
$ws.Range("O377").Value = 'Provincia del Elquí'
$ws.Range("P377").Value = 940

# New row 378 (appended data point).
$ws.Range("A378").Value = 8
$ws.Range("B378").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C378").Value = 'Coquimbo'
$ws.Range("D378").Value = 45007
$ws.Range("D378").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E378").Value = 4
$ws.Range("F378").Value = 100112031
$ws.Range("G378").Value = 'Poroto verde'
$ws.Range("H378").Value = 'Magnum'
$ws.Range("I378").Value = 'Primera'
$ws.Range("J378").Value = 560
$ws.Range("K378").Value = 23500
$ws.Range("L378").Value = 24000
$ws.Range("M378").Value = 23750
$ws.Range("N378").Value = '$/malla 25 kilos'
$ws.Range("O378").Value = 'Provincia de Limarí'
$ws.Range("P378").Value = 950
$ws.Range("Q378").Value = 25
$ws.Range("R378").Value = 'Hortaliza'
